# "Add files via upload" — Spencer Worms adds his "Track Model" risk table
# to the previously-empty sixth worksheet, tweaks a couple of other sheets'
# selections / row heights, and leaves the workbook with the new sheet active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Track Model" sheet: build the risk table that Spencer added.
# ---------------------------------------------------------------------
$wsTM = $wb.Worksheets.Item("Track Model")

# -- Cell values, entered in the same order the shared-string table
#    records them in (body text first, title last) so new strings line
#    up with the target file. --
$wsTM.Range("D3").Value = "Schedule more time to work on module than expected"
$wsTM.Range("A4").Value = "Computer crash, losing files"
$wsTM.Range("D4").Value = "Back up all code to GitHub to ensure that the files still exist"
$wsTM.Range("A5").Value = "Computer loses power while working"
$wsTM.Range("D5").Value = "Save files often to both my computer and the team GitHub"
$wsTM.Range("A6").Value = "Withdraw from class"
$wsTM.Range("A1").Value = "Spencer Worms's Risk Table"
$wsTM.Range("D6").Value = "Work already done on the module will be given to the group and they will have to work on it"
$wsTM.Range("A7").Value = "Personal emergency"
$wsTM.Range("D7").Value = "Work with group to determine if I can work on my module still, or get someone else to help work on it"
$wsTM.Range("A8").Value = "Communication failure with the other modules"
$wsTM.Range("D8").Value = "Work closely with the other team members to ensure good communication and complete work eariler to account of time to fix the issue"
$wsTM.Range("A9").Value = "Have issues projecting my modules when presenting"
$wsTM.Range("D9").Value = "Test the computer before hand and have backup on a flash drive to give to someone else to present on their computer"

$wsTM.Range("A2").Value = "Potential Risk"
$wsTM.Range("B2").Value = "Severity"
$wsTM.Range("C2").Value = "Likelihood"
$wsTM.Range("D2").Value = "Action Plan"

$wsTM.Range("A3").Value = "Failing to complete module on time "
$wsTM.Range("B3").Value = "High"
$wsTM.Range("C3").Value = "Medium"
$wsTM.Range("B4").Value = "High"
$wsTM.Range("C4").Value = "Low"
$wsTM.Range("B5").Value = "Medium"
$wsTM.Range("C5").Value = "Low"
$wsTM.Range("B6").Value = "High"
$wsTM.Range("C6").Value = "Low"
$wsTM.Range("B7").Value = "Medium"
$wsTM.Range("C7").Value = "Low"
$wsTM.Range("B8").Value = "Low"
$wsTM.Range("C8").Value = "Medium"
$wsTM.Range("B9").Value = "Medium"
$wsTM.Range("C9").Value = "Medium"

# -- Column widths for the two text columns. --
$wsTM.Columns.Item(1).ColumnWidth = 39.8
$wsTM.Columns.Item(4).ColumnWidth = 39.8

# -- Row heights for the wrapped body rows (2 lines vs 3 lines tall). --
$wsTM.Rows.Item(3).RowHeight = 30
$wsTM.Rows.Item(4).RowHeight = 30
$wsTM.Rows.Item(5).RowHeight = 30
$wsTM.Rows.Item(6).RowHeight = 30
$wsTM.Rows.Item(7).RowHeight = 30
$wsTM.Rows.Item(8).RowHeight = 45
$wsTM.Rows.Item(9).RowHeight = 45

# -- Formatting: header row gets a thin bottom border, body rows get a
#    full thin box border; both wrap their text. Build each combination
#    once off to the side, then copy/paste the resulting format so the
#    style table only gains the two new cell styles it needs. --
$wsTM.Range("Z1").Borders.Item(9).LineStyle = 1
$wsTM.Range("Z1").WrapText = $true
$wsTM.Range("Z2").Borders.LineStyle = 1
$wsTM.Range("Z2").WrapText = $true

$wsTM.Range("Z1").Copy() | Out-Null
$wsTM.Range("A2:D2").PasteSpecial(-4122) | Out-Null
$wsTM.Range("Z2").Copy() | Out-Null
$wsTM.Range("A3:D9").PasteSpecial(-4122) | Out-Null

$wsTM.Range("Z1:Z2").Clear() | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Selection / active-sheet bookkeeping on the other sheets.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Project RA").Range("D8").Select() | Out-Null
$wb.Worksheets.Item("Train Controller").Range("F25").Select() | Out-Null

$wsCTC = $wb.Worksheets.Item("CTC Office")
$wsCTC.Rows.Item(11).RowHeight = 135

# "Track Model" ends up the active tab, with G8 the remembered selection.
$wsTM.Range("G8").Select() | Out-Null
